# Weekly update: insert two new data rows (for Variedad "Española" and "Madrigal")
# at the top of the data table (rows 19-20), pushing all existing data rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 19 (this shifts rows 19..116 down to 21..118)
$ws.Range("A19:A20").EntireRow.Insert()

# --- New row 19 ---
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "Vega Modelo de Temuco"
$ws.Range("C19").Value = "La Araucanía"
$ws.Range("D19").Value = 44462
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = 100112013
$ws.Range("G19").Value = "Alcachofa"
$ws.Range("H19").Value = "Española"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 12000
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 12000
$ws.Range("N19").Value = "`$/caja 30 unidades"
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 400
$ws.Range("Q19").Value = 30
$ws.Range("R19").Value = "Hortaliza"

# --- New row 20 ---
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 44462
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 100112013
$ws.Range("G20").Value = "Alcachofa"
$ws.Range("H20").Value = "Madrigal"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 12000
$ws.Range("N20").Value = "`$/caja 40 unidades"
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 300
$ws.Range("Q20").Value = 40
$ws.Range("R20").Value = "Hortaliza"
